$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 ("2021年") continues the yearly time series that ends at row 11
# ("2020年"). Copy row 11's formatting down first so the new row-label cell
# (A12) picks up the same bold/centered/bordered style used by every other
# year label in column A, then fill in the row's own values.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A12").Value = "2021年"

# Columns with no reported figure for 2021 stay blank, matching the pattern
# of the other year rows.
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 2485
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 1388
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = 106
$ws.Range("L12").Value = 406
$ws.Range("M12").Value = 2
$ws.Range("N12").Value = 53
$ws.Range("O12").Value = 7
$ws.Range("P12").Value = ""
$ws.Range("Q12").Value = 46
$ws.Range("R12").Value = 250
$ws.Range("S12").Value = ""
$ws.Range("T12").Value = 136
$ws.Range("U12").Value = ""
$ws.Range("V12").Value = ""
$ws.Range("W12").Value = 72
$ws.Range("X12").Value = ""
$ws.Range("Y12").Value = ""
